# Auto-generated edit script: updates market-price-derived columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook
$sheets = $wb.Worksheets

$ws = $sheets.Item(1)
$ws.Range("H9").Value = 1974.4445
$ws.Range("I9").Value = 1974.4445
$ws.Range("J9").Value = 0.0
$ws.Range("K9").Value = 1974.4445
$ws.Range("L9").Value = 0.0
$ws.Range("M9").Value = -1805.4445
$ws.Range("N9").ClearContents()

$ws.Range("H40").Value = 1232.2307
$ws.Range("I40").Value = 1111.039
$ws.Range("J40").Value = 1898.7858
$ws.Range("K40").Value = 1111.039
$ws.Range("L40").Value = 1898.7858
$ws.Range("M40").Value = -936.039
$ws.Range("N40").Value = -2248.7858

$ws.Range("H70").Value = 2535.262
$ws.Range("I70").Value = 2736.4517
$ws.Range("J70").Value = 1968.2727
$ws.Range("K70").Value = 8209.3551
$ws.Range("L70").Value = 5904.8181
$ws.Range("M70").Value = -7939.355100000001
$ws.Range("N70").Value = -6444.8181

$ws.Range("H73").Value = 2535.262
$ws.Range("I73").Value = 2736.4517
$ws.Range("J73").Value = 1968.2727
$ws.Range("K73").Value = 8209.3551
$ws.Range("L73").Value = 5904.8181
$ws.Range("M73").Value = -7273.355100000001
$ws.Range("N73").Value = -7776.8181

$ws.Range("H74").Value = 4463.409
$ws.Range("I74").Value = 3916.5386
$ws.Range("J74").Value = 5253.3335
$ws.Range("K74").Value = 3916.5386
$ws.Range("L74").Value = 5253.3335
$ws.Range("M74").Value = -2980.5386
$ws.Range("N74").Value = -7125.3335

$ws.Range("H77").Value = 4463.409
$ws.Range("I77").Value = 3916.5386
$ws.Range("J77").Value = 5253.3335
$ws.Range("K77").Value = 19582.693
$ws.Range("L77").Value = 26266.6675
$ws.Range("M77").Value = -14902.693
$ws.Range("N77").Value = -35626.6675

$ws.Range("H88").Value = 3704.5833
$ws.Range("I88").Value = 711.7778
$ws.Range("J88").Value = 5500.2666
$ws.Range("K88").Value = 711.7778
$ws.Range("L88").Value = 5500.2666
$ws.Range("M88").Value = -305.7778
$ws.Range("N88").Value = -6312.2666

$ws.Range("H91").Value = 3704.5833
$ws.Range("I91").Value = 711.7778
$ws.Range("J91").Value = 5500.2666
$ws.Range("K91").Value = 711.7778
$ws.Range("L91").Value = 5500.2666
$ws.Range("M91").Value = 692.2222
$ws.Range("N91").Value = -8308.266599999999

$ws.Range("H120").Value = 38755.0
$ws.Range("J120").Value = 38755.0
$ws.Range("L120").Value = 38755.0
$ws.Range("N120").Value = -48431.0

$ws.Range("H132").Value = 1638.0
$ws.Range("I132").Value = 1022.2909
$ws.Range("J132").Value = 5400.6665
$ws.Range("K132").Value = 3066.8727
$ws.Range("L132").Value = 16201.9995
$ws.Range("M132").Value = -536.8726999999999
$ws.Range("N132").Value = -21261.9995

$ws.Range("H136").Value = 28949.5
$ws.Range("J136").Value = 28949.5
$ws.Range("L136").Value = 28949.5
$ws.Range("N136").Value = -39149.5

$ws.Range("H139").Value = 34745.0
$ws.Range("J139").Value = 34745.0
$ws.Range("L139").Value = 34745.0
$ws.Range("N139").Value = -45025.0

$ws = $sheets.Item(2)
$ws.Range("H41").Value = 10000.0
$ws.Range("I41").Value = 0.0
$ws.Range("J41").Value = 10000.0
$ws.Range("K41").Value = 0.0
$ws.Range("L41").Value = 10000.0
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -10828.0

$ws.Range("H122").Value = 3229.1538
$ws.Range("I122").Value = 2967.3333
$ws.Range("J122").Value = 4328.8
$ws.Range("K122").Value = 8901.999899999999
$ws.Range("L122").Value = 12986.4
$ws.Range("M122").Value = -6451.999899999999
$ws.Range("N122").Value = -17886.4

$ws = $sheets.Item(3)
$ws.Range("H64").Value = 405.97726
$ws.Range("I64").Value = 412.43478
$ws.Range("J64").Value = 398.90475
$ws.Range("K64").Value = 412.43478
$ws.Range("L64").Value = 398.90475
$ws.Range("M64").Value = -187.43478
$ws.Range("N64").Value = -848.9047499999999

$ws.Range("H67").Value = 405.97726
$ws.Range("I67").Value = 412.43478
$ws.Range("J67").Value = 398.90475
$ws.Range("K67").Value = 412.43478
$ws.Range("L67").Value = 398.90475
$ws.Range("M67").Value = 367.56522
$ws.Range("N67").Value = -1958.90475

$ws.Range("H94").Value = 1233.8667
$ws.Range("I94").Value = 1212.375
$ws.Range("J94").Value = 1258.4286
$ws.Range("K94").Value = 1212.375
$ws.Range("L94").Value = 1258.4286
$ws.Range("M94").Value = -761.375
$ws.Range("N94").Value = -2160.4286

$ws.Range("H105").Value = 2565777.2
$ws.Range("I105").Value = 1558.9656
$ws.Range("J105").Value = 10002010.0
$ws.Range("K105").Value = 1558.9656
$ws.Range("L105").Value = 10002010.0
$ws.Range("M105").Value = 188.0344
$ws.Range("N105").Value = -10005504.0

$ws = $sheets.Item(5)
$ws.Range("H131").Value = 1316.0834
$ws.Range("J131").Value = 1333.4468
$ws.Range("L131").Value = 4000.3404
$ws.Range("N131").Value = -14080.3404

$ws = $sheets.Item(6)
$ws.Range("H80").Value = 3990.196
$ws.Range("I80").Value = 4624.8237
$ws.Range("J80").Value = 2720.9412
$ws.Range("K80").Value = 4624.8237
$ws.Range("L80").Value = 2720.9412
$ws.Range("M80").Value = -3626.8237
$ws.Range("N80").Value = -4716.9412

$ws.Range("H83").Value = 3990.196
$ws.Range("I83").Value = 4624.8237
$ws.Range("J83").Value = 2720.9412
$ws.Range("K83").Value = 23124.1185
$ws.Range("L83").Value = 13604.706
$ws.Range("M83").Value = -18132.1185
$ws.Range("N83").Value = -23588.706

$ws.Range("H135").Value = 57827.145
$ws.Range("I135").Value = 60000.0
$ws.Range("J135").Value = 54930.0
$ws.Range("K135").Value = 60000.0
$ws.Range("L135").Value = 54930.0
$ws.Range("M135").Value = -54930.0
$ws.Range("N135").Value = -65070.0

$ws = $sheets.Item(7)
$ws.Range("H82").Value = 1697.9412
$ws.Range("I82").Value = 999.5
$ws.Range("J82").Value = 2078.9092
$ws.Range("K82").Value = 999.5
$ws.Range("L82").Value = 2078.9092
$ws.Range("M82").Value = -638.5
$ws.Range("N82").Value = -2800.9092

$ws.Range("H85").Value = 1697.9412
$ws.Range("I85").Value = 999.5
$ws.Range("J85").Value = 2078.9092
$ws.Range("K85").Value = 999.5
$ws.Range("L85").Value = 2078.9092
$ws.Range("M85").Value = 248.5
$ws.Range("N85").Value = -4574.9092

$ws.Range("H93").Value = 1102.0
$ws.Range("I93").Value = 1002.46155
$ws.Range("J93").Value = 1533.3334
$ws.Range("K93").Value = 1002.46155
$ws.Range("L93").Value = 1533.3334
$ws.Range("M93").Value = 245.53845
$ws.Range("N93").Value = -4029.3334

$ws.Range("H132").Value = 13132.904
$ws.Range("I132").Value = 4889.8
$ws.Range("J132").Value = 20626.637
$ws.Range("K132").Value = 14669.4
$ws.Range("L132").Value = 61879.91099999999
$ws.Range("M132").Value = -12139.4
$ws.Range("N132").Value = -66939.911

$ws = $sheets.Item(8)
$ws.Range("H34").Value = 6760.857
$ws.Range("I34").Value = 6326.0
$ws.Range("J34").Value = 6833.3335
$ws.Range("K34").Value = 6326.0
$ws.Range("L34").Value = 6833.3335
$ws.Range("M34").Value = -6123.0
$ws.Range("N34").Value = -7239.3335
